# Commit: "Must clear code and fix dropdown elements and interaction with parent account"
#
# This script updates the SLA Expiration Date cell (AB2) on sheet "Hoja1"
# from a real date value to a plain text value "11/5/2020" (clearing the
# date number format in favor of a text format), renames the built-in
# "Hyperlink" cell style to its Spanish equivalent, and moves the active
# selection to AB2 (matching where the edited cell now lives).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- AB2: "SLA Expiration Date" -------------------------------------------
# Was a date serial (44140 => 11/5/2020) formatted with numFmtId 14 (date).
# Becomes a text cell "11/5/2020" formatted with numFmtId 49 (text/"@").
$cell = $ws.Range("AB2")
$cell.NumberFormat = "@"
$cell.Value = "11/5/2020"

# --- Rename the built-in Hyperlink cell style to the localized name -------
$hyperlinkStyle = $wb.Styles.Item("Hyperlink")
$hyperlinkStyle.Name = "Hipervínculo"

# --- Update the selection / active cell to AB2 -----------------------------
$ws.Range("AB2").Select() | Out-Null
